# "fill out strategic priorities and re-render"
#
# The underlying re-render shrank the template's base body-text size from
# 10pt to 9pt (Normal: sz 20 -> 18). Heading 4 didn't carry its own explicit
# size before (it inherited Normal's 10pt), so it now needs an explicit
# 10pt override to keep rendering the same now that Normal shrank underneath
# it. A handful of other styles (Heading 5, Body Text, Definition Term,
# Definition, Table Caption, Image Caption) already carried an explicit 9pt
# override that is now redundant with the new Normal default, but setting
# them to the same effective 9pt keeps every one of them rendering exactly
# as before.

$d = $word.ActiveDocument

# Normal: 10pt -> 9pt (the template's new base body size).
$d.Styles("Normal").Font.Size = 9

# Heading 4: pin at 10pt explicitly so it keeps its prior rendered size now
# that it no longer inherits a 10pt Normal.
$d.Styles("Heading4").Font.Size = 10

# These already rendered at 9pt (explicit sz=18 override on a then-10pt
# Normal); keep them at 9pt now that it matches the new Normal default.
$d.Styles("Heading5").Font.Size = 9
$d.Styles("BodyText").Font.Size = 9
$d.Styles("DefinitionTerm").Font.Size = 9
$d.Styles("Definition").Font.Size = 9
$d.Styles("TableCaption").Font.Size = 9
$d.Styles("ImageCaption").Font.Size = 9
